$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete() | Out-Null
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("B6").Select() | Out-Null
